$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New student row (row 7)
$ws.Range("A7").Value = "1ru353"
$ws.Range("B7").Value = "Some Name"
$ws.Range("E7").Value = 819273

# F7 (DOB) must stay plain text like the rest of column F ("1987/02/02").
# A direct string assignment gets auto-parsed as a date by the engine,
# so copy the identical existing text cell (F6) instead, which preserves
# the shared-string text type without touching any cell styles.
$ws.Range("F6").Copy()
$ws.Range("F7").PasteSpecial()

$ws.Range("G7").Value = "a"
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = "m"
$ws.Range("J7").Value = "BE"
$ws.Range("K7").Value = "CSE"

$ws.Range("A7").Select()
